# Weekly update: insert two new price rows at the top of the data block
# (row 20), pushing the existing rows down by two. This mirrors the
# "Fruta / hortaliza, semanal" rolling update pattern seen in this
# workbook: the newest entries are inserted near the top and everything
# else shifts down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at row 20 (existing row 20 -> row 22, etc.)
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# New row 20
$ws.Cells.Item(20, 1).Value2 = 10
$ws.Cells.Item(20, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(20, 3).Value2 = "La Araucanía"
$ws.Cells.Item(20, 4).Value2 = 44452
$ws.Cells.Item(20, 5).Value2 = 9
$ws.Cells.Item(20, 6).Value2 = "Fruta"
$ws.Cells.Item(20, 7).Value2 = 100107
$ws.Cells.Item(20, 8).Value2 = "Otros"
$ws.Cells.Item(20, 9).Value2 = 100107002
$ws.Cells.Item(20, 10).Value2 = "Chirimoya"
$ws.Cells.Item(20, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(20, 12).Value2 = "Especial"
$ws.Cells.Item(20, 13).Value2 = 45
$ws.Cells.Item(20, 14).Value2 = 3500
$ws.Cells.Item(20, 15).Value2 = 3500
$ws.Cells.Item(20, 16).Value2 = 3500
$ws.Cells.Item(20, 17).Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(20, 18).Value2 = "Provincia del Elquí"
$ws.Cells.Item(20, 19).Value2 = 3500
$ws.Cells.Item(20, 20).Value2 = 1

# New row 21
$ws.Cells.Item(21, 1).Value2 = 10
$ws.Cells.Item(21, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(21, 3).Value2 = "La Araucanía"
$ws.Cells.Item(21, 4).Value2 = 44452
$ws.Cells.Item(21, 5).Value2 = 9
$ws.Cells.Item(21, 6).Value2 = "Fruta"
$ws.Cells.Item(21, 7).Value2 = 100107
$ws.Cells.Item(21, 8).Value2 = "Otros"
$ws.Cells.Item(21, 9).Value2 = 100107002
$ws.Cells.Item(21, 10).Value2 = "Chirimoya"
$ws.Cells.Item(21, 11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(21, 12).Value2 = "Primera"
$ws.Cells.Item(21, 13).Value2 = 65
$ws.Cells.Item(21, 14).Value2 = 3000
$ws.Cells.Item(21, 15).Value2 = 3000
$ws.Cells.Item(21, 16).Value2 = 3000
$ws.Cells.Item(21, 17).Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(21, 18).Value2 = "Provincia del Elquí"
$ws.Cells.Item(21, 19).Value2 = 3000
$ws.Cells.Item(21, 20).Value2 = 1
